# Apply the "laboratoria 4" update to zadania_dziel_i_rzadz workbook.
#
# Sheet mapping (by tab order in the workbook):
#   1 = zadanie_1     (merge-sort trace table)
#   2 = zadanie_2     (insertion/selection-sort style table that gets filled in)
#   3 = zadanie_2_1   (only a selection/view change)
#   4 = Arkusz1       (hidden helper sheet, untouched)

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# ---------------------------------------------------------------------------
# zadanie_1 (sheet 1): row 28 used to hold a couple of stray values in
# G28/H28 and an empty O28; they get cleared and a full row of "sorted
# so far" values (style copied from the matching row above) is written
# into O28:X28 instead.
# ---------------------------------------------------------------------------
$ws1.Range("G28:H28").ClearContents()

$ws1.Range("T27").Copy() | Out-Null
$ws1.Range("O28:X28").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws1.Range("O28").Value = 1
$ws1.Range("P28").Value = 25
$ws1.Range("Q28").Value = 35
$ws1.Range("R28").Value = 42
$ws1.Range("S28").Value = 55
$ws1.Range("T28").Value = 59
$ws1.Range("U28").Value = 65
$ws1.Range("V28").Value = 68
$ws1.Range("W28").Value = 70
$ws1.Range("X28").Value = 79

# ---------------------------------------------------------------------------
# zadanie_2 (sheet 2): fill in the previously-empty practice table
# (columns I/J/K hold the student's answers, the other columns hold the
# worked trace of the algorithm).
# ---------------------------------------------------------------------------
$ws2.Range("I8").Value = 0
$ws2.Range("J8").Value = 1
$ws2.Range("K8").Value = 2

$ws2.Range("I9").Value = 2
$ws2.Range("J9").Value = 4
$ws2.Range("K9").Value = 3
$ws2.Range("O9").Value = 42
$ws2.Range("P9").Value = 68
$ws2.Range("R9").Value = 35
$ws2.Range("S9").Value = 65
$ws2.Range("T9").Value = 1
$ws2.Range("V9").Value = 25
$ws2.Range("W9").Value = 79
$ws2.Range("Y9").Value = 59
$ws2.Range("Z9").Value = 70
$ws2.Range("AA9").Value = 55

$ws2.Range("I10").Value = 5
$ws2.Range("J10").Value = 6
$ws2.Range("K10").Value = 6

$ws2.Range("I11").Value = 7
$ws2.Range("J11").Value = 9
$ws2.Range("K11").Value = 8

$ws2.Range("I13").Value = 0
$ws2.Range("J13").Value = 0
$ws2.Range("M13").Value = 42
$ws2.Range("O13").Value = 68
$ws2.Range("Q13").Value = 35
$ws2.Range("S13").Value = 65
$ws2.Range("T13").Value = 1
$ws2.Range("V13").Value = 25
$ws2.Range("X13").Value = 79
$ws2.Range("Z13").Value = 59
$ws2.Range("AB13").Value = 70
$ws2.Range("AC13").Value = 55

$ws2.Range("I14").Value = 1
$ws2.Range("J14").Value = 1

$ws2.Range("I16").Value = 2
$ws2.Range("J16").Value = 2
$ws2.Range("R16").Value = 65
$ws2.Range("T16").Value = 1
$ws2.Range("AA16").Value = 70
$ws2.Range("AC16").Value = 55

$ws2.Range("I17").Value = 3
$ws2.Range("J17").Value = 4
$ws2.Range("K17").Value = 4

$ws2.Range("I19").Value = 3
$ws2.Range("J19").Value = 3

$ws2.Range("I20").Value = 4
$ws2.Range("J20").Value = 4

$ws2.Range("I22").Value = 5
$ws2.Range("J22").Value = 5

$ws2.Range("I23").Value = 6
$ws2.Range("J23").Value = 6

$ws2.Range("I25").Value = 7
$ws2.Range("J25").Value = 7

$ws2.Range("I26").Value = 8
$ws2.Range("J26").Value = 9
$ws2.Range("K26").Value = 9

$ws2.Range("I28").Value = 8
$ws2.Range("J28").Value = 8

$ws2.Range("I29").Value = 9

# ---------------------------------------------------------------------------
# View state: update each sheet's frozen-pane / selection to match where
# the author was last working, then leave zadanie_2 as the active tab.
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("M29").Select()

$ws3.Activate()
$ws3.Range("E31").Select()

$ws2.Activate()
$ws2.Range("O16").Select()
